# IO-fixes: FxE matrix improvements.
#
# A new data row (Country/Entity/Parameter/Type/Flow/Value = CHE / ext_biomass /
# output / configuration_fxe / biomass / 1) is inserted above the existing
# row 7 ("output_efficiency" / "constant_fxe" / biomass / 1) of the FxE
# (flow x entity) matrix. This pushes every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7 (shifts rows 7..end down by one row).
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new "output" / "configuration_fxe" entry.
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "ext_biomass"
$ws.Range("C7").Value = "output"
$ws.Range("D7").Value = "configuration_fxe"
$ws.Range("F7").Value = "biomass"
$ws.Range("G7").Value = 1
$ws.Range("K7").Style = "Hyperlink"
$ws.Range("K7").ClearContents()

# The row insertion does not automatically re-anchor the two existing
# hyperlinks (originally on K8 and K30); recreate them at their shifted
# locations (K9 and K31) pointing at the same DOI target.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K31"), "https://doi.org/10.1016/j.esr.2019.100379", "", "Persistent link using digital object identifier") | Out-Null
$ws.Hyperlinks.Add($ws.Range("K9"), "https://doi.org/10.1016/j.esr.2019.100379", "", "Persistent link using digital object identifier") | Out-Null
$ws.Range("K31").Style = "Hyperlink"
$ws.Range("K9").Style = "Hyperlink"

# Likewise, re-apply the autofilter so its range grows from A5:L572 to A5:L573.
$ws.AutoFilterMode = $false
$ws.Range("A5:L573").AutoFilter()

# And update the _FilterDatabase defined name to match the new filter range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$5:`$L`$573"
    }
}

# Match the post-edit selection state recorded in the workbook.
$ws.Range("E7").Select()
